$wb = $excel.ActiveWorkbook

# Rename the "settings" sheet to "Settings" (creation of the settings window).
$ws = $wb.Worksheets.Item("settings")
$ws.Name = "Settings"

# Update the selection on the Settings sheet to E28.
$ws.Range("E28").Select()

